# T-P_POO-NeoDarbellay.xlsx — "Finished the project" update
#
# Adds a 3rd entry block (2 new rows) to the "Jour 7" section of the
# Journal sheet (rows 62-64 after the edit), describing the final
# wrap-up work (more waves, testing, releasing to GitHub), and shifts
# everything below down by two rows accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- 1. Insert two new rows right before the (old) row 62 ---------------
# This pushes the "end of Jour 7 block" (old row 62), its totals row
# (old row 63) and everything after it down by two rows, and Excel
# automatically widens/shifts the dependent merged cells, the SUM()
# formulas and the plain data-validation sqrefs that span the insertion
# point.
$ws.Rows("62:63").Insert()

# --- 2. Give the two new rows the same look as the row above them -------
# (borders / number formats / alignment) instead of the blank default
# formatting Insert() leaves behind.
$ws.Range("A61:G61").Copy()
$ws.Range("A62:G63").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows("62:63").RowHeight = 18

# --- 3. Fill in the three "Jour 7" rows that now make up the tail of ----
#        the block: 62 (new), 63 (new) and 64 (was row 62, now holds the
#        last entry of the day instead of being blank).
$ws.Range("A62").Value = "Coding"
$ws.Range("C62").Value = 25
$ws.Range("D62").Value = "Creation of more waves"
$ws.Range("E62").Value = "Finished"
$ws.Range("F62").Value2 = 0.70833333333333337

$ws.Range("D63").Value = "Testing the game out"
$ws.Range("A63").Value = "Test"
$ws.Range("C63").Value = 5
$ws.Range("E63").Value = "Finished"
$ws.Range("F63").Value2 = 0.71180555555555547

$ws.Range("D64").Value = "Releasing the final version to GitHub"
$ws.Range("A64").Value = "GitHub"
$ws.Range("C64").Value = 10
$ws.Range("E64").Value = "Finished"
$ws.Range("F64").NumberFormat = "h:mm"
$ws.Range("F64").Value2 = 0.71875

# --- 4. Update the print area to match the grown sheet ------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$G`$72"

$wb.Save()
